# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values computed from the regenerated save_data (replaces old Strike#-based values)
$kValues = @(0, 1, 0, 0, 1, 0, 1, 2, 1, 1, 1)

$row = 2
foreach ($k in $kValues) {
    $ws.Cells.Item($row, 7).Value = $k
    $row++
}
